$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Paaiy259"
$ws.Range("B2").Value = 23091477
$ws.Range("C2").Value = "qvmifem19"
$ws.Range("D2").Value = "Qy9k2B$%"
$ws.Range("F2").Value = "lICnwatu"
$ws.Range("G2").Value = "BEos"
